$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (M2:T2)
$ws.Range("M2").Value = 1.321358333333333
$ws.Range("N2").Value = 3.964075
$ws.Range("O2").Value = 0.06904315418552966
$ws.Range("P2").Value = 0.06904315418552966
$ws.Range("Q2").Value = 0.2254642533222222
$ws.Range("R2").Value = 2.0291782799
$ws.Range("S2").Value = 0.06904315418552966
$ws.Range("T2").Value = 0.06904315418552966

# Row 3 updates (O3, P3, S3, T3)
$ws.Range("O3").Value = 0.4558096119837698
$ws.Range("P3").Value = 0.4558096119837698
$ws.Range("S3").Value = 0.4558096119837698
$ws.Range("T3").Value = 0.4558096119837698

# Row 4 updates (M4:T4)
$ws.Range("M4").Value = 9.093439666666667
$ws.Range("N4").Value = 27.280319
$ws.Range("O4").Value = 0.4751472338307006
$ws.Range("P4").Value = 0.4751472338307005
$ws.Range("Q4").Value = 1.551619672616445
$ws.Range("R4").Value = 13.964577053548
$ws.Range("S4").Value = 0.4751472338307006
$ws.Range("T4").Value = 0.4751472338307005
